$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2789.95
$ws.Range("I43").Value = 2919.9333
$ws.Range("K43").Value = 2919.9333
$ws.Range("M43").Value = -2850.9333
$ws.Range("H64").Value = 93527.55
$ws.Range("I64").Value = 1000000
$ws.Range("J64").Value = 2880.3
$ws.Range("K64").Value = 1000000
$ws.Range("L64").Value = 2880.3
$ws.Range("M64").Value = -999752
$ws.Range("N64").Value = -3376.3
$ws.Range("H67").Value = 93527.55
$ws.Range("I67").Value = 1000000
$ws.Range("J67").Value = 2880.3
$ws.Range("K67").Value = 1000000
$ws.Range("L67").Value = 2880.3
$ws.Range("M67").Value = -999142
$ws.Range("N67").Value = -4596.3
$ws.Range("H98").Value = 34055.082
$ws.Range("I98").Value = 1088.9286
$ws.Range("J98").Value = 80207.7
$ws.Range("K98").Value = 1088.9286
$ws.Range("L98").Value = 80207.7
$ws.Range("M98").Value = 409.0714
$ws.Range("N98").Value = -83203.7
$ws.Range("H122").Value = 34055.082
$ws.Range("I122").Value = 1088.9286
$ws.Range("J122").Value = 80207.7
$ws.Range("K122").Value = 3266.7858
$ws.Range("L122").Value = 240623.1
$ws.Range("M122").Value = -816.7857999999997
$ws.Range("N122").Value = -245523.1
$ws.Range("H132").Value = 16346.564
$ws.Range("I132").Value = 2392.283
$ws.Range("J132").Value = 98521.78
$ws.Range("K132").Value = 7176.849
$ws.Range("L132").Value = 295565.34
$ws.Range("M132").Value = -4646.849
$ws.Range("N132").Value = -300625.34
$ws.Range("H138").Value = 1745.5522
$ws.Range("I138").Value = 1506.3871
$ws.Range("J138").Value = 1951.5
$ws.Range("K138").Value = 4519.1613
$ws.Range("L138").Value = 5854.5
$ws.Range("M138").Value = 620.8387000000002
$ws.Range("N138").Value = -16134.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3403.6572
$ws.Range("I61").Value = 2138.1428
$ws.Range("J61").Value = 4247.3335
$ws.Range("K61").Value = 2138.1428
$ws.Range("L61").Value = 4247.3335
$ws.Range("M61").Value = -1926.1428
$ws.Range("N61").Value = -4671.3335
$ws.Range("H74").Value = 1135.6316
$ws.Range("I74").Value = 639.3333
$ws.Range("K74").Value = 639.3333
$ws.Range("M74").Value = 234.6667
$ws.Range("H77").Value = 1135.6316
$ws.Range("I77").Value = 639.3333
$ws.Range("K77").Value = 3196.6665
$ws.Range("M77").Value = 1171.3335
$ws.Range("H132").Value = 1382.3954
$ws.Range("I132").Value = 1023.5278
$ws.Range("J132").Value = 3228
$ws.Range("K132").Value = 3070.5834
$ws.Range("L132").Value = 9684
$ws.Range("M132").Value = -540.5834
$ws.Range("N132").Value = -14744
$ws.Range("H136").Value = 3403.6572
$ws.Range("I136").Value = 2138.1428
$ws.Range("J136").Value = 4247.3335
$ws.Range("K136").Value = 6414.428400000001
$ws.Range("L136").Value = 12742.0005
$ws.Range("M136").Value = -3864.428400000001
$ws.Range("N136").Value = -17842.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2064.1924
$ws.Range("I105").Value = 1779.6666
$ws.Range("K105").Value = 1779.6666
$ws.Range("M105").Value = -32.66660000000002
$ws.Range("H134").Value = 3143.33
$ws.Range("I134").Value = 1321.0605
$ws.Range("J134").Value = 4040.8657
$ws.Range("K134").Value = 3963.1815
$ws.Range("L134").Value = 12122.5971
$ws.Range("M134").Value = -1428.1815
$ws.Range("N134").Value = -17192.5971

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 43000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 43000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 43000
$ws.Range("N55").Value = -43630
$ws.Range("M55").Value = ""
$ws.Range("H99").Value = 1688.2
$ws.Range("I99").Value = 1720.2222
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 1720.2222
$ws.Range("L99").Value = 1400
$ws.Range("M99").Value = -222.2221999999999
$ws.Range("N99").Value = -4396
$ws.Range("H122").Value = 80760
$ws.Range("I122").Value = 100766.664
$ws.Range("J122").Value = 733.3333
$ws.Range("K122").Value = 302299.992
$ws.Range("L122").Value = 2199.9999
$ws.Range("M122").Value = -299849.992
$ws.Range("N122").Value = -7099.9999
$ws.Range("H126").Value = 1688.2
$ws.Range("I126").Value = 1720.2222
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 5160.6666
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -2690.6666
$ws.Range("N126").Value = -9140
$ws.Range("H132").Value = 50106.414
$ws.Range("I132").Value = 1434.3914
$ws.Range("J132").Value = 236682.5
$ws.Range("K132").Value = 4303.174199999999
$ws.Range("L132").Value = 710047.5
$ws.Range("M132").Value = -1773.174199999999
$ws.Range("N132").Value = -715107.5
$ws.Range("H134").Value = 453115.94
$ws.Range("I134").Value = 969.5417
$ws.Range("J134").Value = 2003332.1
$ws.Range("K134").Value = 2908.6251
$ws.Range("L134").Value = 6009996.300000001
$ws.Range("M134").Value = -373.6251000000002
$ws.Range("N134").Value = -6015066.300000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2970.75
$ws.Range("J52").Value = 2970.75
$ws.Range("L52").Value = 8912.25
$ws.Range("N52").Value = -9444.25
$ws.Range("H98").Value = 739.8125
$ws.Range("I98").Value = 521.1667
$ws.Range("J98").Value = 871
$ws.Range("K98").Value = 1563.5001
$ws.Range("L98").Value = 2613
$ws.Range("M98").Value = -65.50009999999997
$ws.Range("N98").Value = -5609
$ws.Range("H113").Value = 3869.1936
$ws.Range("I113").Value = 5635.05
$ws.Range("J113").Value = 658.5454999999999
$ws.Range("K113").Value = 16905.15
$ws.Range("L113").Value = 1975.6365
$ws.Range("M113").Value = -14735.15
$ws.Range("N113").Value = -6315.6365
$ws.Range("H121").Value = 310168.9
$ws.Range("J121").Value = 775224.75
$ws.Range("L121").Value = 2325674.25
$ws.Range("N121").Value = -2328294.25
$ws.Range("H131").Value = 890.5700000000001
$ws.Range("J131").Value = 906.70526
$ws.Range("L131").Value = 2720.11578
$ws.Range("N131").Value = -12800.11578

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 25414.285
$ws.Range("J95").Value = 25414.285
$ws.Range("L95").Value = 25414.285
$ws.Range("N95").Value = -30906.285

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2627
$ws.Range("I100").Value = 2287.3333
$ws.Range("J100").Value = 2966.6667
$ws.Range("K100").Value = 2287.3333
$ws.Range("L100").Value = 2966.6667
$ws.Range("M100").Value = -1746.3333
$ws.Range("N100").Value = -4048.6667
$ws.Range("H132").Value = 4268.5557
$ws.Range("I132").Value = 2439.2
$ws.Range("J132").Value = 5344.647
$ws.Range("K132").Value = 7317.599999999999
$ws.Range("L132").Value = 16033.941
$ws.Range("M132").Value = -4787.599999999999
$ws.Range("N132").Value = -21093.941
$ws.Range("H136").Value = 2105.3901
$ws.Range("I136").Value = 1610.9
$ws.Range("J136").Value = 3454
$ws.Range("K136").Value = 4832.700000000001
$ws.Range("L136").Value = 10362
$ws.Range("M136").Value = -2282.700000000001
$ws.Range("N136").Value = -15462

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 17394
$ws.Range("J40").Value = 17394
$ws.Range("L40").Value = 17394
$ws.Range("N40").Value = -17692
$ws.Range("H80").Value = 33197.777
$ws.Range("J80").Value = 33197.777
$ws.Range("L80").Value = 33197.777
$ws.Range("N80").Value = -35193.777
$ws.Range("H83").Value = 33197.777
$ws.Range("J83").Value = 33197.777
$ws.Range("L83").Value = 99593.33100000001
$ws.Range("N83").Value = -109577.331
$ws.Range("H103").Value = 37956.285
$ws.Range("J103").Value = 37956.285
$ws.Range("L103").Value = 37956.285
$ws.Range("N103").Value = -40300.285
$ws.Range("H122").Value = 2041374.2
$ws.Range("I122").Value = 2381436.8
$ws.Range("K122").Value = 7144310.399999999
$ws.Range("M122").Value = -7141860.399999999
$ws.Range("H126").Value = 2943796.5
$ws.Range("I126").Value = 3270551.5
$ws.Range("K126").Value = 9811654.5
$ws.Range("M126").Value = -9809184.5
$ws.Range("H132").Value = 2254.1292
$ws.Range("I132").Value = 1782.8948
$ws.Range("J132").Value = 3000.25
$ws.Range("K132").Value = 5348.6844
$ws.Range("L132").Value = 9000.75
$ws.Range("M132").Value = -2818.6844
$ws.Range("N132").Value = -14060.75
$ws.Range("H136").Value = 23802.756
$ws.Range("I136").Value = 56410.223
$ws.Range("K136").Value = 169230.669
$ws.Range("M136").Value = -166680.669
